# Commit: add template example (new "vlinder_metadata_template" worksheet),
# based on the existing "mocca_template" sheet.

$wb = $excel.ActiveWorkbook

$srcSheet = $wb.Worksheets.Item("mocca_template")
$vlinderWs = $wb.Worksheets.Item("vlinder_template")

# --- create the new sheet as a copy of mocca_template, placed right after it ---
$srcSheet.Copy($null, $srcSheet)
$ws = $wb.Worksheets.Item(3)
$ws.Name = "vlinder_metadata_template"

# --- row 2 ---
$ws.Range("A2").Value = "_ID"
$ws.Range("B2").Value = "ID"

# --- row 3 ---
$ws.Range("A3").Value = "name"
$ws.Range("B3").Value = "VLINDER"
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = "object"

# --- row 4 ---
$ws.Range("A4").Value = "lat"
$ws.Range("B4").Value = "lat"
$ws.Range("F4").ClearContents()

# --- row 5 ---
$ws.Range("A5").Value = "lon"
$ws.Range("B5").Value = "lon"
$ws.Range("F5").ClearContents()

# --- row 6 ---
$ws.Range("A6").Value = "location"
$ws.Range("B6").Value = "stad"
$ws.Range("F6").ClearContents()

# --- row 7 ---
$ws.Range("A7").Value = "call_name"
$ws.Range("B7").Value = "benaming"
$ws.Range("E2").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = "object"

# --- row 8 ---
$ws.Range("A8").Value = "network"
$ws.Range("B8").Value = "Network"
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").Value = "object"

# --- rows 9-17: blank filler rows use the same style as the other data cells ---
$ws.Range("B2").Copy()
$ws.Range("B9:E17").PasteSpecial(-4122)

# --- rows 9-17: wipe any leftover sample content copied from the mocca template ---
$ws.Range("A9:F17").ClearContents()

$excel.CutCopyMode = 0

# --- page setup / header & footer for the new metadata template sheet ---
$ws.PageSetup.LeftMargin = 56.7
$ws.PageSetup.RightMargin = 56.7
$ws.PageSetup.TopMargin = 75.8
$ws.PageSetup.BottomMargin = 75.8
$ws.PageSetup.HeaderMargin = 56.7
$ws.PageSetup.FooterMargin = 56.7
$ws.PageSetup.CenterHeader = "&""Times New Roman,Regular""&12&A"
$ws.PageSetup.CenterFooter = "&""Times New Roman,Regular""&12Page &P"

# --- refresh the zoom level on all three template sheets ---
$vlinderWs.Activate()
$vlinderWs.Application.ActiveWindow.Zoom = 200
$vlinderWs.Range("A1").Select() | Out-Null

$srcSheet.Activate()
$srcSheet.Application.ActiveWindow.Zoom = 200

$ws.Activate()
$ws.Application.ActiveWindow.Zoom = 200
$ws.Range("E13").Select() | Out-Null
